$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing data down
$ws.Rows.Item(2).Insert()

# Fix formatting of the newly inserted row to match the rest of the data rows
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2:D2").ClearFormats()
$ws.Range("E2").Clear()

# Write out the recalculated forecast data for every row
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 2.070003986395053
$ws.Range("D2").Value = 2008

$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("D3").Value = 2009

$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 0.517569958955022
$ws.Range("D4").Value = 2009

$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("D5").Value = 2010

$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = -3.956152295564896
$ws.Range("D6").Value = 2010

$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 0.7825601129312298
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = -0.2328395085068102

$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = 1.234995474941392
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 0.7327527981400461

$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 1.171373351779592
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 0.6458049593451864

$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 0.899360810820804
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 0.4113848771853501

$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 0.4797371259343874
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 0.5448268972618964

$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 0.9010266119894084
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 1.203907967581297

$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = 0.1088602047940146
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 0.1671551101610103

$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = 0.02019328874804938
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = -0.03860754389363175

$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = -0.075394216261504
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = -0.06463514052835739

$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 0.1729981757035093
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 0.2629870913912535

$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = -0.07548837955325682
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 0.05928147027902675

$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 0.09752710595589686
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 0.1987429576382649

$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = -0.5735475396625112
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 0.02406984837131088

$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = -0.5280591151586633
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 0.1903092973221776

$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 0.02883110668334687
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 0.241498802789164

$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 0.07201851318385799
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 0.2843016498274009

$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 0.8258453722611359
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 0.288064297781454

$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 0.7252300059688022
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 0.2052430644269299

$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 0.4640111827386662
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = -0.1620498231152179

$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 0.3727661260635617
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = -0.8612142616933327

$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = -0.9065026814729205
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = -0.1365403697986656

$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = -0.6243248145489155
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 0.07475225043114264

$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = -0.7901161779547028
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = -0.1949185821441768

$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = -0.801759526476209
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 0.06491682578968483

$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = -0.2930109800340586
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = -0.1644047560850792

$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -0.6491730431770759
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = -0.4370777949570193

$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = -1.503583188367719
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = -0.5877299932847579

$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = -1.103489789942047
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = 1.323658311025055

$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = 2.632698787096288
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = 0.4893323826990148

$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = 1.311904119834839
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = -0.1159018519404809

$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = 1.067534122491809
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = -0.112565850764601

$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 0.9704846793491928
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = -0.6989646400249128

$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = -0.3568974718008655
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = 0.04689880979749095

$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = -1.338216592160768
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = -0.2487719682984557

$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = -0.9795431199870586
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = -0.3324688493351879

$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = -0.7009264669202708
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = -0.0234350458557242

$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = -0.05370673382950608
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = -0.03417477517112522

$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = 0.2512652100014945
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 0.041441321352087

$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = 0.1829021030556488
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 0.01245506629512505

$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = 0.3928252664241905
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = 0.196134499498668

$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = 0.2464401331885524
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = -0.02784774425726999

$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = 0.9693451788297391
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 0.08029846083614789

$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = 0.6979546684258597
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = -0.04101132194430646

$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = 0.3224026462283813
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = -0.6671574593505647

$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = -2.205730080079726
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = -0.2015885781823656

$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = -1.551451534890558
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = -0.1079783528070921

$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = -2.11737366557071
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = -0.2835476113072333
